$d = $word.ActiveDocument

# Correct the parcellation count: 64 -> 68 regions of interest.
$d.Content.Find.Execute("segmented into 64 regions", $true, $false, $false, $false, $false,
                         $true, 1, $false, "segmented into 68 regions", 2)
